$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; existing rows 20-43 shift down to 21-44.
$ws.Rows(20).Insert()

# Populate the newly inserted row 20 with the latest week's data.
$ws.Cells.Item(20, 1).Value = 6
$ws.Cells.Item(20, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(20, 3).Value = "Metropolitana"
$ws.Cells.Item(20, 4).Value = 44803
$ws.Cells.Item(20, 5).Value = 13
$ws.Cells.Item(20, 6).Value = 100112035
$ws.Cells.Item(20, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 400
$ws.Cells.Item(20, 11).Value = 16000
$ws.Cells.Item(20, 12).Value = 18000
$ws.Cells.Item(20, 13).Value = 16850
$ws.Cells.Item(20, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(20, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(20, 16).Value = 1123
$ws.Cells.Item(20, 17).Value = 15
$ws.Cells.Item(20, 18).Value = "Hortaliza"
